# Add the 0610 "blue tank" titration row (new temp probe) to CRMAccuracyData.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# --- New row 41 data ---
$ws.Range("A41").Value = 20210610
$ws.Range("B41").Value = 2229.3150000000001
$ws.Range("C41").Value = 2224.4699999999998
$ws.Range("D41").Formula = "=100*(B41-C41)/C41"
$ws.Range("E41").Value = 180
$ws.Range("F41").Value = "CRM OPENED 20210526"

# Re-assert D40's formula so it keeps evaluating correctly alongside the new row.
$ws.Range("D40").Formula = "=100*(B40-C40)/C40"

# Move the active selection to A42, just below the newly added row, as on save.
$ws.Range("A42").Select() | Out-Null
